# Update EPEX spot price workbook with the latest day of data:
#  - "Prix Spot": add a new date column Z (09-jul) with its 24 hourly prices
#  - "Gaz": append a new row (2025-07-07 / 33.4)
#  - "CO2": append a new row (2025-07-07 / 69.95999999999999)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column Z = 09-jul
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, styled like the other day headers (bold / bordered / centered)
$wsPrix.Range("Z1").Value = "09-jul"
$wsPrix.Range("Y1").Copy()
$wsPrix.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

$hourlyPrices = @(
    71.98999999999999,
    39.37,
    33.59,
    23.37,
    15.28,
    54.77,
    15,
    61.47,
    69.13,
    61.61,
    20.99,
    0,
    11.68,
    1.38,
    0,
    5,
    2.93,
    30.02,
    64.28,
    103.96,
    111.23,
    95.16,
    115.91,
    97.37
)

for ($i = 0; $i -lt $hourlyPrices.Count; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 26).Value = $hourlyPrices[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append 2025-07-07 / 33.4
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Format the new date cell as text first so the ISO-looking string isn't
# auto-converted to a serial date, then strip the formatting back off (the
# source cells in column A carry no explicit style) while keeping the text.
$wsGaz.Range("A23").NumberFormat = "@"
$wsGaz.Range("A23").Value = "2025-07-07"
$wsGaz.Range("A23").ClearFormats()

$wsGaz.Range("B23").Value = 33.4

# ---------------------------------------------------------------------------
# Sheet "CO2": append 2025-07-07 / 69.95999999999999
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A23").NumberFormat = "@"
$wsCo2.Range("A23").Value = "2025-07-07"
$wsCo2.Range("A23").ClearFormats()

$wsCo2.Range("B23").Value = 69.95999999999999
